# The page's picture (the single InlineShape sitting in the first
# paragraph) is being dropped from the document, leaving the paragraph
# that hosted it empty - the second (already empty) paragraph and the
# section properties are untouched.
$d = $word.ActiveDocument

while ($d.InlineShapes.Count -gt 0) {
    $d.InlineShapes(1).Delete()
}

Write-Output ("InlineShapes remaining: " + $d.InlineShapes.Count)
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
